# Updates "Price" (D) and "Volume(1h)" (E) columns for the cryptos list.
# D-column values that parse as plain numbers are entered with a leading
# apostrophe so Excel keeps them as text (matching the original inlineStr
# cells, e.g. "1.002"), then ClearFormats() strips the resulting
# quote-prefix formatting so the cell style stays at its original default.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.463.84'
$ws.Range("E2").Value = '  -1.16%  '
$ws.Range("D3").Value = '2.091.82'
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("E5").Value = '  -1.04%  '
$ws.Range("E6").Value = '  +0.18%  '
$ws.Range("E7").Value = '  -0.44%  '
$ws.Range("D8").Value = '''0.4381'
$ws.Range("E8").Value = '  -0.75%  '
$ws.Range("D9").Value = '''54.04'
$ws.Range("E9").Value = '  +15.65%  '
$ws.Range("E10").Value = '  -1.61%  '
$ws.Range("D11").Value = '''1.153'
$ws.Range("E11").Value = '  -2.88%  '
$ws.Range("D12").Value = '''24.27'
$ws.Range("E12").Value = '  -4.50%  '
$ws.Range("D13").Value = '2.093.28'
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("D14").Value = '''6.689'
$ws.Range("E14").Value = '  -1.32%  '
$ws.Range("D15").Value = '''7.675'
$ws.Range("E15").Value = '  -2.73%  '
$ws.Range("D16").Value = '''95.92'
$ws.Range("E16").Value = '  -2.52%  '
$ws.Range("D17").Value = '''1.002'
$ws.Range("E17").Value = '  +0.19%  '
$ws.Range("D18").Value = '''0.00001123'
$ws.Range("E18").Value = '  -1.53%  '
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").Value = '''19.18'
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("D22").Value = '''6.260'
$ws.Range("E22").Value = '  -2.36%  '
$ws.Range("D23").Value = '30.496.35'
$ws.Range("E23").Value = '  -1.43%  '
$ws.Range("D24").Value = '''12.28'
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("D25").Value = '''2.321'
$ws.Range("E25").Value = '  +2.89%  '
$ws.Range("D26").Value = '2.332.87'
$ws.Range("E26").Value = '  -1.55%  '
$ws.Range("D27").Value = '''22.23'
$ws.Range("E27").Value = '  -3.36%  '
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("D29").Value = '''163.39'
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("D30").Value = '''131.55'
$ws.Range("E30").Value = '  -1.80%  '
$ws.Range("D31").Value = '''1.185'
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("E32").Value = '  -0.18%  '
$ws.Range("D33").Value = '''1.659'
$ws.Range("E33").Value = '  +5.47%  '
$ws.Range("D34").Value = '''6.163'
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").Value = '''3.899'
$ws.Range("E35").Value = '  -3.61%  '
$ws.Range("D36").Value = '''10.02'
$ws.Range("E36").Value = '  +4.41%  '
$ws.Range("D37").Value = '''0.02565'
$ws.Range("E37").Value = '  -1.45%  '
$ws.Range("D38").Value = '''0.06826'
$ws.Range("E38").Value = '  +0.89%  '
$ws.Range("D39").Value = '''5.479'
$ws.Range("E39").Value = '  -2.11%  '
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("D41").Value = '''0.2253'
$ws.Range("E41").Value = '  -0.86%  '
$ws.Range("D42").Value = '''0.6876'
$ws.Range("E42").Value = '  +0.31%  '
$ws.Range("D43").Value = '''1.252'
$ws.Range("E43").Value = '  -0.75%  '
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").Value = '''13.91'
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("D46").Value = '''0.6335'
$ws.Range("E46").Value = '  -1.49%  '
$ws.Range("E47").Value = '  -3.12%  '
$ws.Range("D48").Value = '''3.625'
$ws.Range("E48").Value = '  -1.52%  '
$ws.Range("D49").Value = '''1.230'
$ws.Range("E49").Value = '  +6.36%  '
$ws.Range("E50").Value = '  -3.37%  '
$ws.Range("D51").Value = '''81.78'
$ws.Range("E51").Value = '  -1.74%  '

$ws.Range("D2:D51").ClearFormats()
